# Agregando la funcion que corrije detalles de comparacion - unificando las
# funciones en un solo bloque de codigo.
#
# Consolida las filas de la tabla: cada combinacion producto/cantidad queda
# en una unica fila con todos los precios completos (PlazaVea, Tottus,
# Metro), eliminando las filas que tenian datos incompletos ("-") o que ya
# no son necesarias tras la unificacion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos finales consolidados (filas 2 a 5 de la hoja).
$data = @(
    @("AZUCAR RUBIA//AZUCAR RUBIA//AZUCAR RUBIA", "DULFINA", "5", "KG", "BOLSA", "21.90", "22.9", "26.45"),
    @("AZUCAR RUBIA//AZUCAR RUBIA//AZUCAR RUBIA", "DULFINA", "1", "KG", "BOLSA", "4.89", "5.2", "5.70"),
    @("AZUCAR RUBIA//AZUCAR RUBIA", "DULFINA", "5", "KG", "BOLSA", "21.90", "22.9", "26.45"),
    @("AZUCAR RUBIA//AZUCAR RUBIA", "DULFINA", "1", "KG", "BOLSA", "4.89", "5.2", "5.70")
)

# Elimina las filas de datos anteriores (2 a 11), dejando unicamente el
# encabezado de la tabla.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -ge 2) {
    $ws.Range("A2:H" + $lastRow).EntireRow.Delete() | Out-Null
}

# Escribe las filas consolidadas a partir de la fila 2. Los valores que
# parecen numericos ("5", "21.90", etc.) deben seguir almacenados como texto
# (tal como en el resto del archivo), asi que se escriben primero como una
# formula de texto y luego se convierten a valores estaticos mediante
# copiar/pegado especial, evitando que Excel los reinterprete como numeros.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $col = $j + 1
        $text = $data[$i][$j]
        $cell = $ws.Cells.Item($row, $col)
        if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
            $escaped = $text.Replace('"', '""')
            $cell.Formula = '="' + $escaped + '"'
            $cell.Copy() | Out-Null
            $cell.PasteSpecial(-4163) | Out-Null
        } else {
            $cell.Value = $text
        }
    }
}

$excel.CutCopyMode = 0
